$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(1).Insert()
$ws.Range("E1").Value = "Target"
$ws.Range("F1").Value = "Total"
$ws.Range("G3").Value = 68.6

$cws = $wb.Worksheets.Item("Chart1")
$co = $cws.ChartObjects().Item(1)
$chart = $co.Chart
$ser1 = $chart.SeriesCollection(1)
$ser1.Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$13,Sheet1!`$E`$2:`$E`$13,1)"
$ser1.Name = "Target"
$ser2 = $chart.SeriesCollection(2)
$ser2.Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$13,Sheet1!`$F`$2:`$F`$13,2)"
$ser2.Name = "Actual"
